$wb = $excel.ActiveWorkbook

# The "Slovakia" sheet is the template the new "Italy" sheet is modeled on;
# insert the new sheet right after it (it becomes the active tab).
$slovakia = $wb.Worksheets.Item("Slovakia")
$italy = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $slovakia)
$italy.Name = "Italy"

# Declare the merges up front so that copying formats afterwards does not
# need to synthesize new border/merge styling (keeps cellXfs identical to
# the template instead of growing styles.xml).
$italy.Range("A1:D1").MergeCells = $true
$italy.Range("C2:D2").MergeCells = $true

# Column widths, matching the template.
$italy.Columns("A").ColumnWidth = $slovakia.Columns("A").ColumnWidth
$italy.Columns("B").ColumnWidth = $slovakia.Columns("B").ColumnWidth
$italy.Columns("C").ColumnWidth = $slovakia.Columns("C").ColumnWidth
$italy.Columns("D").ColumnWidth = $slovakia.Columns("D").ColumnWidth

# Formats (per-cell styles, merges) for the populated block.
$slovakia.Range("A1:D5").Copy()
$italy.Range("A1:D5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$slovakia.Range("A6:A9").Copy()
$italy.Range("A6:A9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# The "Transmission Unit / ... / ... Keysafe" rows (Slovakia rows 10-12) are
# skipped for Italy; rows 13-15 (Black Box, Wg, Attached Functionality)
# become Italy's rows 10-12.
$slovakia.Range("A13:A15").Copy()
$italy.Range("A10:A12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Values: reuse the shared text from Slovakia where identical, and set the
# two genuinely new strings for Italy.
$italy.Range("A1").Value = $slovakia.Range("A1").Value()

$italy.Range("A2").Value = $slovakia.Range("A2").Value()
$italy.Range("B2").Value = "Italy Market"
$italy.Range("C2").Value = $slovakia.Range("C2").Value()

$italy.Range("A3").Value = $slovakia.Range("A3").Value()
$italy.Range("D3").Value = $slovakia.Range("D3").Value()

$italy.Range("A4").Value = $slovakia.Range("A4").Value()
$italy.Range("B4").Value = "NGC-3145/T2157"
$italy.Range("B4").ClearFormats()
$italy.Range("D4").Value = $slovakia.Range("D4").Value()

$italy.Range("D5").Value = $slovakia.Range("D5").Value()

$italy.Range("A7").Value = $slovakia.Range("A7").Value()
$italy.Range("A8").Value = $slovakia.Range("A8").Value()
$italy.Range("A9").Value = $slovakia.Range("A9").Value()

$italy.Range("A10").Value = $slovakia.Range("A13").Value()
$italy.Range("A11").Value = $slovakia.Range("A14").Value()
$italy.Range("A12").Value = $slovakia.Range("A15").Value()

# Match the recorded view state: Italy ends up active with A17 selected
# (the row just past the used range), Slovakia loses its tab selection.
$italy.Range("A17").Select()
